$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 / Row 5: the SP25 course shown in row 4 and row 5 were swapped ---
$ws.Range("D4").Value = "CPSC 6127"
$ws.Range("E4").Value = "Contemporary Issues in Database Management Systems"

$ws.Range("D5").Value = "CPSC 6179"
$ws.Range("E5").Value = "Software Project Planning and Management"

# --- Row 6: move "CPSC 6000" / "Graduate Exit Examination..." from B/C to D/E ---
$ws.Range("B6").Copy()
$ws.Range("D6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D6").Value = "CPSC 6000"

$ws.Range("C6").Copy()
$ws.Range("E6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E6").Value = "Graduate Exit Examination in Computer Science"

$ws.Range("B6").Clear()
$ws.Range("C6").Clear()

# --- Row 8: swap course counts ---
$ws.Range("B8").Value = "Courses: 2"
$ws.Range("D8").Value = "Courses: 3"

# --- Formatting: taller title row, wrap text on every centered/styled cell ---
$ws.Rows.Item(2).RowHeight = 30

$ws.Range("B2").WrapText = $true

$ws.Range("B3").WrapText = $true
$ws.Range("D3").WrapText = $true
$ws.Range("F3").WrapText = $true

$ws.Range("B4").WrapText = $true
$ws.Range("C4").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Range("E4").WrapText = $true

$ws.Range("B5").WrapText = $true
$ws.Range("C5").WrapText = $true
$ws.Range("D5").WrapText = $true
$ws.Range("E5").WrapText = $true

$ws.Range("D6").WrapText = $true
$ws.Range("E6").WrapText = $true

$ws.Range("B8").WrapText = $true
$ws.Range("D8").WrapText = $true
$ws.Range("F8").WrapText = $true
